$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(1.071367240641678,0.9781698075374834,0.9211656525764624,0.8979920336477107,0.8941474819982034,0.9208528970079612,1.039187697688931,1.272958446809639,1.445742302250551,1.524568878466027,1.554450490150714,1.548013553945282,1.527026626157181,1.514175630075897,1.440595281546621,1.395513480550051,1.369604974067556,1.360836492529188,1.40031031147339,1.533190141968475,1.620219173278542,1.573753545803356,1.398141634182934,1.209535015349729)
$CValues = @(0.1885825876654792,0.166363432426607,0.1527054039895575,0.1471360221240445,0.1462110188490158,0.1526303077058344,0.1809247515872698,0.2362803560961879,0.2768650382289195,0.2953086650025796,0.3022899565873445,0.3007865437387522,0.2958830797426515,0.2928791804480966,0.275659309758737,0.2650905519347475,0.2590099502237706,0.2569508730049961,0.2662157952518101,0.2973234268498857,0.3176369944344799,0.3067969014800269,0.2657070866954996,0.2213197265586189)
$DValues = @(0.5889528017668937,0.5856100738657943,0.5838469299353193,0.5832011884097739,0.5830983578977822,0.5838379266520803,0.5877401846873909,0.597688945145876,0.6064013327122098,0.610670274831989,0.6123307957948327,0.6119712172665288,0.6108060056679676,0.6100980059771928,0.6061284987209774,0.6037716251868517,0.6024447805022817,0.6020004739391425,0.6040195410525655,0.6111470632253315,0.6160615822973625,0.6134151566033381,0.6039073705979376,0.5947513809599485)
$EValues = @(0.2177370710720936,0.217636365364541,0.217679123566743,0.2177228433125116,0.2177316909855094,0.2176796067237099,0.2176806359680938,0.2185130602574787,0.2196320120476081,0.2202515015224087,0.2205019880550303,0.2204473339530253,0.2202717904755787,0.2201663359944561,0.2195937510882331,0.2192707929466629,0.2190954333456041,0.2190378449546593,0.2193040962062796,0.2203229202376669,0.2210814605815159,0.2206681285324947,0.2192890076788103,0.2181988513590731)
$GValues = @(0.002488314989996804,0.002492099847007265,0.002494547020349213,0.002495575355536721,0.00249574799076503,0.002494560762687218,0.002489594486411195,0.002480829105640048,0.002474976313349098,0.002472439869744557,0.002471497403758135,0.002471699580180404,0.002472361971675305,0.002472770050703137,0.002475144603982027,0.002476633527306169,0.00247750178340499,0.002477797801308732,0.002476473801404307,0.002472166922614139,0.002469457180993329,0.002470893838834125,0.002476545975222285,0.002483096811976609)
$JValues = @(0.09830449300914879,0.09890801516050729,0.09932231810599035,0.09950214340054231,0.09953266709584518,0.09932469878670247,0.09850350669833219,0.09724044319465364,0.0965245889488564,0.09624508105306262,0.09614588080188824,0.09616694974567253,0.09623678659573187,0.09628042910250301,0.09654378460304258,0.09671716785308604,0.09682123548550692,0.09685721650579637,0.09669826144813598,0.0962160934513534,0.0959396953470133,0.0960836677962682,0.09670679536819193,0.09754491457637826)
$MValues = @(0.495408168708245,0.4701469265253024,0.4548234435517671,0.448626223810983,0.4476000382981766,0.4547396742071186,0.4866593832413955,0.5507330078283132,0.5987089246288946,0.6207303526277457,0.6290975253055393,0.6272942564898187,0.6214181621121355,0.6178225433747855,0.5972737196321987,0.5847179907484161,0.577514815448211,0.575079135377905,0.5860526512415305,0.62314335131326,0.6475481833985697,0.6345079264209659,0.5854492036240089,0.533241134061889)
$OValues = @(4.783484108390468,4.797194863521185,4.808889663886561,4.814477426378261,4.815454865675946,4.808961696458226,4.787530462396319,4.771587689327788,4.775904517333288,4.781377375217346,4.783956587479167,4.783378535320196,4.781579399983116,4.780543436750065,4.77561764945699,4.773496136606525,4.772606063448706,4.772361340626333,4.773687788831268,4.782094079129934,4.790542863828421,4.785762511782877,4.773600116255153,4.773094908849657)

$colMap = @{
    'B' = $BValues
    'C' = $CValues
    'D' = $DValues
    'E' = $EValues
    'G' = $GValues
    'J' = $JValues
    'M' = $MValues
    'O' = $OValues
}

foreach ($col in $colMap.Keys) {
    $values = $colMap[$col]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Range("$col$row").Value = $values[$i]
    }
}
